# "reshaped gui and moved dup to load_data"
#
# Adds a "Breed" column (column B) next to the existing "Name" column in
# the members sheet: Scooby Doo is a dog, everyone else (Shaggy, Velma,
# Fred, and Daphne twice) is human.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the data rows first -- this is the order new values land in the
# shared-strings table -- then the header last.
$ws.Range("B2").Value = "dog"
$ws.Range("B3").Value = "human"
$ws.Range("B4").Value = "human"
$ws.Range("B5").Value = "human"
$ws.Range("B6").Value = "human"
$ws.Range("B7").Value = "human"

$ws.Range("B1").Value = "Breed"
$ws.Range("B1").Font.Bold = $true

# Leave the selection where the author's last click landed.
$ws.Range("D7").Select() | Out-Null
